# chore: update Sheets via scheduled runner
# Refreshes cached market-price-derived figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) for a handful of leve rows across several
# crafting-job sheets in the Carbuncle_Profits workbook.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 121
$ws.Cells.Item(121, 8).Value = 1117.9166
$ws.Cells.Item(121, 10).Value = 1180
$ws.Cells.Item(121, 12).Value = 3540
$ws.Cells.Item(121, 14).Value = -7034
# Row 132
$ws.Cells.Item(132, 8).Value = 1093.027
$ws.Cells.Item(132, 9).Value = 1045.0883
$ws.Cells.Item(132, 11).Value = 3135.2649
$ws.Cells.Item(132, 13).Value = -605.2648999999997
# Row 137
$ws.Cells.Item(137, 8).Value = 295858.9
$ws.Cells.Item(137, 9).Value = 435878.12
$ws.Cells.Item(137, 10).Value = 3091.4546
$ws.Cells.Item(137, 11).Value = 1307634.36
$ws.Cells.Item(137, 12).Value = 9274.363799999999
$ws.Cells.Item(137, 13).Value = -1305084.36
$ws.Cells.Item(137, 14).Value = -14374.3638
# Row 138
$ws.Cells.Item(138, 8).Value = 3610.7407
$ws.Cells.Item(138, 9).Value = 1291.45
$ws.Cells.Item(138, 10).Value = 10237.286
$ws.Cells.Item(138, 11).Value = 3874.35
$ws.Cells.Item(138, 12).Value = 30711.858
$ws.Cells.Item(138, 13).Value = 1265.65
$ws.Cells.Item(138, 14).Value = -40991.858

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 1421.7222
$ws.Cells.Item(61, 9).Value = 879.875
$ws.Cells.Item(61, 10).Value = 5756.5
$ws.Cells.Item(61, 11).Value = 879.875
$ws.Cells.Item(61, 12).Value = 5756.5
$ws.Cells.Item(61, 13).Value = -667.875
$ws.Cells.Item(61, 14).Value = -6180.5
# Row 132
$ws.Cells.Item(132, 8).Value = 2426.4092
$ws.Cells.Item(132, 9).Value = 1201.5652
$ws.Cells.Item(132, 10).Value = 3767.9048
$ws.Cells.Item(132, 11).Value = 3604.6956
$ws.Cells.Item(132, 12).Value = 11303.7144
$ws.Cells.Item(132, 13).Value = -1074.6956
$ws.Cells.Item(132, 14).Value = -16363.7144
# Row 134
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 14).ClearContents()
# Row 136
$ws.Cells.Item(136, 8).Value = 1421.7222
$ws.Cells.Item(136, 9).Value = 879.875
$ws.Cells.Item(136, 10).Value = 5756.5
$ws.Cells.Item(136, 11).Value = 2639.625
$ws.Cells.Item(136, 12).Value = 17269.5
$ws.Cells.Item(136, 13).Value = -89.625
$ws.Cells.Item(136, 14).Value = -22369.5
# Row 137
$ws.Cells.Item(137, 8).Value = 53580
$ws.Cells.Item(137, 10).Value = 53580
$ws.Cells.Item(137, 12).Value = 53580
$ws.Cells.Item(137, 14).Value = -63780

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Cells.Item(94, 8).Value = 502
$ws.Cells.Item(94, 9).Value = 359.7143
$ws.Cells.Item(94, 10).Value = 1000
$ws.Cells.Item(94, 11).Value = 359.7143
$ws.Cells.Item(94, 12).Value = 1000
$ws.Cells.Item(94, 13).Value = 91.28570000000002
$ws.Cells.Item(94, 14).Value = -1902
# Row 99
$ws.Cells.Item(99, 8).Value = 1810
$ws.Cells.Item(99, 9).Value = 1177.5
$ws.Cells.Item(99, 10).Value = 2500
$ws.Cells.Item(99, 11).Value = 1177.5
$ws.Cells.Item(99, 12).Value = 2500
$ws.Cells.Item(99, 13).Value = 320.5
$ws.Cells.Item(99, 14).Value = -5496
# Row 134
$ws.Cells.Item(134, 8).Value = 1710.8611
$ws.Cells.Item(134, 9).Value = 1429.5
$ws.Cells.Item(134, 10).Value = 2153
$ws.Cells.Item(134, 11).Value = 4288.5
$ws.Cells.Item(134, 12).Value = 6459
$ws.Cells.Item(134, 13).Value = -1753.5
$ws.Cells.Item(134, 14).Value = -11529

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 2475
$ws.Cells.Item(31, 9).Value = 2124.9092
$ws.Cells.Item(31, 10).Value = 2672.487
$ws.Cells.Item(31, 11).Value = 2124.9092
$ws.Cells.Item(31, 12).Value = 2672.487
$ws.Cells.Item(31, 13).Value = -1829.9092
$ws.Cells.Item(31, 14).Value = -3262.487
# Row 34
$ws.Cells.Item(34, 8).Value = 2475
$ws.Cells.Item(34, 9).Value = 2124.9092
$ws.Cells.Item(34, 10).Value = 2672.487
$ws.Cells.Item(34, 11).Value = 2124.9092
$ws.Cells.Item(34, 12).Value = 2672.487
$ws.Cells.Item(34, 13).Value = -1922.9092
$ws.Cells.Item(34, 14).Value = -3076.487
# Row 94
$ws.Cells.Item(94, 8).Value = 1856.8572
$ws.Cells.Item(94, 9).Value = 2632
$ws.Cells.Item(94, 10).Value = 1426.2222
$ws.Cells.Item(94, 11).Value = 2632
$ws.Cells.Item(94, 12).Value = 1426.2222
$ws.Cells.Item(94, 13).Value = -2181
$ws.Cells.Item(94, 14).Value = -2328.2222
# Row 132
$ws.Cells.Item(132, 8).Value = 2020.3256
$ws.Cells.Item(132, 9).Value = 857.9286
$ws.Cells.Item(132, 10).Value = 4190.1333
$ws.Cells.Item(132, 11).Value = 2573.7858
$ws.Cells.Item(132, 12).Value = 12570.3999
$ws.Cells.Item(132, 13).Value = -43.78579999999965
$ws.Cells.Item(132, 14).Value = -17630.3999
# Row 134
$ws.Cells.Item(134, 8).Value = 2574.7083
$ws.Cells.Item(134, 9).Value = 2336.4285
$ws.Cells.Item(134, 10).Value = 4242.6665
$ws.Cells.Item(134, 11).Value = 7009.2855
$ws.Cells.Item(134, 12).Value = 12727.9995
$ws.Cells.Item(134, 13).Value = -4474.2855
$ws.Cells.Item(134, 14).Value = -17797.9995
# Row 140
$ws.Cells.Item(140, 8).Value = 47954.855
$ws.Cells.Item(140, 10).Value = 47954.855
$ws.Cells.Item(140, 12).Value = 47954.855
$ws.Cells.Item(140, 14).Value = -58314.855

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 115
$ws.Cells.Item(115, 8).Value = 2877.4443
$ws.Cells.Item(115, 9).Value = 0
$ws.Cells.Item(115, 10).Value = 2877.4443
$ws.Cells.Item(115, 11).Value = 0
$ws.Cells.Item(115, 12).Value = 8632.332900000001
$ws.Cells.Item(115, 13).ClearContents()
$ws.Cells.Item(115, 14).Value = -10982.3329
# Row 123
$ws.Cells.Item(123, 8).Value = 1403.8667
$ws.Cells.Item(123, 9).Value = 1011.6
$ws.Cells.Item(123, 10).Value = 1600
$ws.Cells.Item(123, 11).Value = 3034.8
$ws.Cells.Item(123, 12).Value = 4800
$ws.Cells.Item(123, 13).Value = -584.8000000000002
$ws.Cells.Item(123, 14).Value = -9700
# Row 125
$ws.Cells.Item(125, 8).Value = 1479.9333
$ws.Cells.Item(125, 9).Value = 799.5
$ws.Cells.Item(125, 10).Value = 1584.6154
$ws.Cells.Item(125, 11).Value = 2398.5
$ws.Cells.Item(125, 12).Value = 4753.8462
$ws.Cells.Item(125, 13).Value = 2521.5
$ws.Cells.Item(125, 14).Value = -14593.8462
# Row 131
$ws.Cells.Item(131, 8).Value = 5264.5835
$ws.Cells.Item(131, 9).Value = 410
$ws.Cells.Item(131, 10).Value = 8732.143
$ws.Cells.Item(131, 11).Value = 1230
$ws.Cells.Item(131, 12).Value = 26196.429
$ws.Cells.Item(131, 13).Value = 3810
$ws.Cells.Item(131, 14).Value = -36276.429

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 5209.49
$ws.Cells.Item(70, 9).Value = 4564.8125
$ws.Cells.Item(70, 10).Value = 6295.263
$ws.Cells.Item(70, 11).Value = 4564.8125
$ws.Cells.Item(70, 12).Value = 6295.263
$ws.Cells.Item(70, 13).Value = -4294.8125
$ws.Cells.Item(70, 14).Value = -6835.263
# Row 73
$ws.Cells.Item(73, 8).Value = 5209.49
$ws.Cells.Item(73, 9).Value = 4564.8125
$ws.Cells.Item(73, 10).Value = 6295.263
$ws.Cells.Item(73, 11).Value = 4564.8125
$ws.Cells.Item(73, 12).Value = 6295.263
$ws.Cells.Item(73, 13).Value = -3628.8125
$ws.Cells.Item(73, 14).Value = -8167.263
# Row 80
$ws.Cells.Item(80, 8).Value = 4955.25
$ws.Cells.Item(80, 9).Value = 5085
$ws.Cells.Item(80, 10).Value = 4773.6
$ws.Cells.Item(80, 11).Value = 5085
$ws.Cells.Item(80, 12).Value = 4773.6
$ws.Cells.Item(80, 13).Value = -4087
$ws.Cells.Item(80, 14).Value = -6769.6
# Row 83
$ws.Cells.Item(83, 8).Value = 4955.25
$ws.Cells.Item(83, 9).Value = 5085
$ws.Cells.Item(83, 10).Value = 4773.6
$ws.Cells.Item(83, 11).Value = 25425
$ws.Cells.Item(83, 12).Value = 23868
$ws.Cells.Item(83, 13).Value = -20433
$ws.Cells.Item(83, 14).Value = -33852
# Row 126
$ws.Cells.Item(126, 8).Value = 1830.975
$ws.Cells.Item(126, 9).Value = 1547.6207
$ws.Cells.Item(126, 10).Value = 2578
$ws.Cells.Item(126, 11).Value = 4642.8621
$ws.Cells.Item(126, 12).Value = 7734
$ws.Cells.Item(126, 13).Value = -2172.8621
$ws.Cells.Item(126, 14).Value = -12674
# Row 132
$ws.Cells.Item(132, 8).Value = 4217
$ws.Cells.Item(132, 9).Value = 4599.7
$ws.Cells.Item(132, 11).Value = 13799.1
$ws.Cells.Item(132, 13).Value = -11269.1

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value = 5833.3335
$ws.Cells.Item(62, 10).Value = 5833.3335
$ws.Cells.Item(62, 12).Value = 5833.3335
$ws.Cells.Item(62, 14).Value = -7081.3335
# Row 65
$ws.Cells.Item(65, 8).Value = 5833.3335
$ws.Cells.Item(65, 10).Value = 5833.3335
$ws.Cells.Item(65, 12).Value = 29166.6675
$ws.Cells.Item(65, 14).Value = -35406.6675
# Row 100
$ws.Cells.Item(100, 8).Value = 378.8
$ws.Cells.Item(100, 9).Value = 378.8
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 11).Value = 757.6
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 13).Value = -216.6
$ws.Cells.Item(100, 14).ClearContents()
